# Update NATMI TPM-derived values in Fam3c-Lifr.xlsx (Sheet1)
# The underlying TPM input changed, so the ligand/receptor expression
# statistics and the derived specificity/edge-weight scores were
# recomputed. All affected cells hold literal numbers (no formulas),
# so we just overwrite them with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 6.086760000000001
$ws.Range("H2").Value = 18.26028
$ws.Range("I2").Value = 0.2669498455274402
$ws.Range("J2").Value = 0.2669498455274402
$ws.Range("M2").Value = 24.91851366666667
$ws.Range("N2").Value = 74.75554099999999
$ws.Range("O2").Value = 0.2924799159147552
$ws.Range("P2").Value = 0.2924799159147553
$ws.Range("Q2").Value = 151.67301224572
$ws.Range("R2").Value = 1365.05711021148
$ws.Range("S2").Value = 0.07807746837332259
$ws.Range("T2").Value = 0.07807746837332261

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 6.086760000000001
$ws.Range("H3").Value = 18.26028
$ws.Range("I3").Value = 0.2669498455274402
$ws.Range("J3").Value = 0.2669498455274402
$ws.Range("O3").Value = 0.4753125595076708
$ws.Range("P3").Value = 0.4753125595076708
$ws.Range("Q3").Value = 246.48560033012
$ws.Range("R3").Value = 2218.37040297108
$ws.Range("S3").Value = 0.1268846143378249
$ws.Range("T3").Value = 0.1268846143378249

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = 6.086760000000001
$ws.Range("H4").Value = 18.26028
$ws.Range("I4").Value = 0.2669498455274402
$ws.Range("J4").Value = 0.2669498455274402
$ws.Range("M4").Value = 19.78346566666667
$ws.Range("N4").Value = 59.350397
$ws.Range("O4").Value = 0.232207524577574
$ws.Range("P4").Value = 0.232207524577574
$ws.Range("Q4").Value = 120.41720748124
$ws.Range("R4").Value = 1083.75486733116
$ws.Range("S4").Value = 0.06198776281629265
$ws.Range("T4").Value = 0.06198776281629265

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.5502435543747958
$ws.Range("J5").Value = 0.5502435543747957
$ws.Range("M5").Value = 24.91851366666667
$ws.Range("N5").Value = 74.75554099999999
$ws.Range("O5").Value = 0.2924799159147552
$ws.Range("P5").Value = 0.2924799159147553
$ws.Range("Q5").Value = 312.6321245697752
$ws.Range("R5").Value = 2813.689121127976
$ws.Range("S5").Value = 0.1609351885161763
$ws.Range("T5").Value = 0.1609351885161763

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.5502435543747958
$ws.Range("J6").Value = 0.5502435543747957
$ws.Range("O6").Value = 0.4753125595076708
$ws.Range("P6").Value = 0.4753125595076708
$ws.Range("R6").Value = 4572.559362373496
$ws.Range("S6").Value = 0.2615376721824824
$ws.Range("T6").Value = 0.2615376721824824

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.5502435543747958
$ws.Range("J7").Value = 0.5502435543747957
$ws.Range("M7").Value = 19.78346566666667
$ws.Range("N7").Value = 59.350397
$ws.Range("O7").Value = 0.232207524577574
$ws.Range("P7").Value = 0.232207524577574
$ws.Range("Q7").Value = 248.2068948998658
$ws.Range("R7").Value = 2233.862054098792
$ws.Range("S7").Value = 0.1277706936761371
$ws.Range("T7").Value = 0.1277706936761371

# Row 8 (MuSCs -> ECs)
$ws.Range("G8").Value = 4.168198333333334
$ws.Range("H8").Value = 12.504595
$ws.Range("I8").Value = 0.1828066000977641
$ws.Range("J8").Value = 0.1828066000977641
$ws.Range("M8").Value = 24.91851366666667
$ws.Range("N8").Value = 74.75554099999999
$ws.Range("O8").Value = 0.2924799159147552
$ws.Range("P8").Value = 0.2924799159147553
$ws.Range("Q8").Value = 103.8653071345439
$ws.Range("R8").Value = 934.7877642108949
$ws.Range("S8").Value = 0.05346725902525633
$ws.Range("T8").Value = 0.05346725902525633

# Row 9 (MuSCs -> FAPs)
$ws.Range("G9").Value = 4.168198333333334
$ws.Range("H9").Value = 12.504595
$ws.Range("I9").Value = 0.1828066000977641
$ws.Range("J9").Value = 0.1828066000977641
$ws.Range("O9").Value = 0.4753125595076708
$ws.Range("P9").Value = 0.4753125595076708
$ws.Range("Q9").Value = 168.7927351311161
$ws.Range("R9").Value = 1519.134616180045
$ws.Range("S9").Value = 0.08689027298736349
$ws.Range("T9").Value = 0.08689027298736347

# Row 10 (MuSCs -> MuSCs)
$ws.Range("G10").Value = 4.168198333333334
$ws.Range("H10").Value = 12.504595
$ws.Range("I10").Value = 0.1828066000977641
$ws.Range("J10").Value = 0.1828066000977641
$ws.Range("M10").Value = 19.78346566666667
$ws.Range("N10").Value = 59.350397
$ws.Range("O10").Value = 0.232207524577574
$ws.Range("P10").Value = 0.232207524577574
$ws.Range("Q10").Value = 82.46140861935723
$ws.Range("R10").Value = 742.152677574215
$ws.Range("S10").Value = 0.04244906808514431
$ws.Range("T10").Value = 0.04244906808514431
